# Adapt column header formatting to respective input file names:
#   "<header>_old" -> "<header>_FV2410"
#   "<header>_new" -> "<header>_FV2504"
# then (re-)expose the header row as an Excel Table (so the new names are
# picked up as the table's column headers) and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header cells (row 1) ------------------------------------
# Columns A:J carry the "_old" suffixed headers, column K is the untouched
# "diff" column, and columns L:U carry the "_new" suffixed headers.
for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $current = [string]$cell.Value
    if ($current.EndsWith("_old")) {
        $cell.Value = $current.Substring(0, $current.Length - 4) + "_FV2410"
    } elseif ($current.EndsWith("_new")) {
        $cell.Value = $current.Substring(0, $current.Length - 4) + "_FV2504"
    }
}

# --- 2. Freeze the header row ------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Turn the used range into an Excel Table, headers taken from row 1 --
$usedRange = $ws.Range("A1:U60")
$tbl = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $usedRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$tbl.Name = "Table1"
